$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Cells.Item(2, 4).Value = '30.063.66'
$ws.Cells.Item(2, 5).Value = '  -0.40%  '

$ws.Cells.Item(3, 4).Value = '1.871.41'
$ws.Cells.Item(3, 5).Value = '  -2.94%  '

Set-TextCell 4 4 '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.31%  '

Set-TextCell 5 4 '319.68'
$ws.Cells.Item(5, 5).Value = '  -1.04%  '

Set-TextCell 6 4 '1.003'
$ws.Cells.Item(6, 5).Value = '  +0.30%  '

Set-TextCell 7 4 '0.5078'
$ws.Cells.Item(7, 5).Value = '  -1.62%  '

Set-TextCell 8 4 '0.3946'
$ws.Cells.Item(8, 5).Value = '  -1.17%  '

$ws.Cells.Item(9, 5).Value = '  -3.34%  '

Set-TextCell 10 4 '42.16'
$ws.Cells.Item(10, 5).Value = '  -2.10%  '

Set-TextCell 11 4 '1.094'
$ws.Cells.Item(11, 5).Value = '  -2.77%  '

Set-TextCell 12 4 '22.69'
$ws.Cells.Item(12, 5).Value = '  +7.12%  '

$ws.Cells.Item(13, 4).Value = '1.870.23'
$ws.Cells.Item(13, 5).Value = '  -2.65%  '

Set-TextCell 14 4 '6.264'
$ws.Cells.Item(14, 5).Value = '  -1.13%  '

Set-TextCell 15 4 '7.165'
$ws.Cells.Item(15, 5).Value = '  -3.02%  '

Set-TextCell 16 4 '1.003'
$ws.Cells.Item(16, 5).Value = '  +0.25%  '

Set-TextCell 17 4 '92.29'
$ws.Cells.Item(17, 5).Value = '  -2.36%  '

Set-TextCell 18 4 '0.00001081'
$ws.Cells.Item(18, 5).Value = '  -3.34%  '

Set-TextCell 19 4 '0.06333'
$ws.Cells.Item(19, 5).Value = '  -6.26%  '

Set-TextCell 20 4 '17.84'
$ws.Cells.Item(20, 5).Value = '  -0.95%  '

Set-TextCell 21 4 '1.003'
$ws.Cells.Item(21, 5).Value = '  +0.35%  '

$ws.Cells.Item(22, 4).Value = '30.037.44'
$ws.Cells.Item(22, 5).Value = '  -0.52%  '

Set-TextCell 23 4 '5.807'
$ws.Cells.Item(23, 5).Value = '  -4.48%  '

Set-TextCell 24 4 '11.05'
$ws.Cells.Item(24, 5).Value = '  -1.75%  '

Set-TextCell 25 4 '2.209'
$ws.Cells.Item(25, 5).Value = '  -0.12%  '

$ws.Cells.Item(26, 4).Value = '2.089.17'
$ws.Cells.Item(26, 5).Value = '  -2.47%  '

Set-TextCell 27 4 '161.40'
$ws.Cells.Item(27, 5).Value = '  +1.16%  '

Set-TextCell 28 4 '21.00'
$ws.Cells.Item(28, 5).Value = '  -0.19%  '

Set-TextCell 29 4 '2.259'
$ws.Cells.Item(29, 5).Value = '  -8.80%  '

Set-TextCell 30 4 '126.80'
$ws.Cells.Item(30, 5).Value = '  -1.83%  '

$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 31 4 '1.046'
$ws.Cells.Item(31, 5).Value = '  -3.22%  '

$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 32 4 '0.1035'
$ws.Cells.Item(32, 5).Value = '  -2.28%  '

Set-TextCell 33 4 '5.894'
$ws.Cells.Item(33, 5).Value = '  -3.51%  '

Set-TextCell 34 4 '3.741'
$ws.Cells.Item(34, 5).Value = '  +2.01%  '

Set-TextCell 35 4 '0.02434'
$ws.Cells.Item(35, 5).Value = '  -2.71%  '

Set-TextCell 36 4 '5.230'
$ws.Cells.Item(36, 5).Value = '  +0.52%  '

Set-TextCell 37 4 '0.06355'
$ws.Cells.Item(37, 5).Value = '  -4.27%  '

Set-TextCell 38 4 '0.2144'
$ws.Cells.Item(38, 5).Value = '  -3.16%  '

$ws.Cells.Item(39, 2).Value = 'ARBITRUM'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 39 4 '1.172'
$ws.Cells.Item(39, 5).Value = '  -5.95%  '

$ws.Cells.Item(40, 2).Value = 'FraxShare'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 40 4 '8.538'
$ws.Cells.Item(40, 5).Value = '  -5.92%  '

Set-TextCell 41 4 '0.6278'
$ws.Cells.Item(41, 5).Value = '  -4.12%  '

Set-TextCell 42 4 '1.209'
$ws.Cells.Item(42, 5).Value = '  -2.83%  '

Set-TextCell 43 4 '11.29'
$ws.Cells.Item(43, 5).Value = '  -1.14%  '

$ws.Cells.Item(44, 2).Value = 'Decentraland'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 44 4 '0.5901'
$ws.Cells.Item(44, 5).Value = '  -4.16%  '

$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 45 4 '12.89'
$ws.Cells.Item(45, 5).Value = '  -2.62%  '

$ws.Cells.Item(46, 2).Value = 'PancakeSwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 46 4 '3.644'
$ws.Cells.Item(46, 5).Value = '  -2.22%  '

$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 47 4 '1.992'
$ws.Cells.Item(47, 5).Value = '  -3.35%  '

$ws.Cells.Item(48, 2).Value = 'EOS'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell 48 4 '1.207'
$ws.Cells.Item(48, 5).Value = '  -2.92%  '

$ws.Cells.Item(49, 2).Value = 'Quant'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 49 4 '121.90'
$ws.Cells.Item(49, 5).Value = '  -2.93%  '

$ws.Cells.Item(50, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 50 4 '1.113'
$ws.Cells.Item(50, 5).Value = '  -3.23%  '

$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 51 4 '76.79'
$ws.Cells.Item(51, 5).Value = '  -3.33%  '
